$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2's COUNTA formula result implicitly recalculates; we just set the
# source data and let Excel recompute the COUNTA(A3:A39) formula in A2.

# Row 5: rename the TextureTag value (B5) from Enemy_Elite_Attack -> Enemy_Elite_SwordSwing
$ws.Range("B5").Value = "Enemy_Elite_SwordSwing"

# New rows 6-8
$ws.Range("A6").Value = "Stab"
$ws.Range("B6").Value = "Enemy_Elite_SwordStab"
$ws.Range("C6").Value = 0.06

$ws.Range("A7").Value = "Up"
$ws.Range("B7").Value = "Enemy_Elite_Up"
$ws.Range("C7").Value = 0.04

$ws.Range("A8").Value = "Death"
$ws.Range("B8").Value = "Enemy_Elite_Death"
$ws.Range("C8").Value = 0.06

$ws.Range("C8").Select()
